$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka1")

$ws.Range("A21").Value = "Donate"
$ws.Range("B21").Value = "Paypal donate funkció beépítése webView segítségével"
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 4
$ws.Range("H21").Value = "Pictori"

$ws.Range("D23").Select()
